# Fruta / hortaliza, semanal
# Update rows 2-10 (D, M, N, O, P, Q, R, S, T) to reflect the weekly reshuffle
# of price-report records described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44174
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 3700
$ws.Range("O2").Value = 3800
$ws.Range("P2").Value = 3747
$ws.Range("Q2").Value = "$/bandeja 2 kilos"
$ws.Range("R2").Value = "Provincia de Linares"
$ws.Range("S2").Value = 1874
$ws.Range("T2").Value = 2

# Row 3
$ws.Range("D3").Value = 44540
$ws.Range("M3").Value = 240
$ws.Range("N3").Value = 3500
$ws.Range("O3").Value = 3800
$ws.Range("P3").Value = 3650
$ws.Range("Q3").Value = "$/bandeja 2 kilos"
$ws.Range("R3").Value = "Región del Maule"
$ws.Range("S3").Value = 1825
$ws.Range("T3").Value = 2

# Row 4
$ws.Range("D4").Value = 44594
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 2500
$ws.Range("O4").Value = 2800
$ws.Range("P4").Value = 2650
$ws.Range("Q4").Value = "$/bandeja 2 kilos"
$ws.Range("R4").Value = "Provincia de Linares"
$ws.Range("S4").Value = 1325
$ws.Range("T4").Value = 2

# Row 5
$ws.Range("D5").Value = 44539
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 3800
$ws.Range("O5").Value = 4000
$ws.Range("P5").Value = 3900
$ws.Range("Q5").Value = "$/bandeja 2 kilos"
$ws.Range("R5").Value = "Región del Maule"
$ws.Range("S5").Value = 1950
$ws.Range("T5").Value = 2

# Row 6
$ws.Range("D6").Value = 44181
$ws.Range("M6").Value = 65
$ws.Range("N6").Value = 3600
$ws.Range("O6").Value = 3800
$ws.Range("P6").Value = 3692
$ws.Range("Q6").Value = "$/bandeja 2 kilos"
$ws.Range("R6").Value = "Provincia de Diguillín"
$ws.Range("S6").Value = 1846
$ws.Range("T6").Value = 2

# Row 7
$ws.Range("D7").Value = 44181
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 1800
$ws.Range("O7").Value = 2000
$ws.Range("P7").Value = 1875
$ws.Range("Q7").Value = "$/envase 1 kilo"
$ws.Range("R7").Value = "Provincia de Diguillín"
$ws.Range("S7").Value = 1875
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("D8").Value = 44596
$ws.Range("M8").Value = 120
$ws.Range("N8").Value = 2500
$ws.Range("O8").Value = 2700
$ws.Range("P8").Value = 2600
$ws.Range("Q8").Value = "$/bandeja 2 kilos"
$ws.Range("R8").Value = "Provincia de Linares"
$ws.Range("S8").Value = 1300
$ws.Range("T8").Value = 2

# Row 9
$ws.Range("D9").Value = 44187
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 2800
$ws.Range("O9").Value = 3000
$ws.Range("P9").Value = 2900
$ws.Range("Q9").Value = "$/bandeja 2 kilos"
$ws.Range("R9").Value = "Provincia de Linares"
$ws.Range("S9").Value = 1450
$ws.Range("T9").Value = 2

# Row 10
$ws.Range("D10").Value = 44187
$ws.Range("M10").Value = 65
$ws.Range("N10").Value = 1400
$ws.Range("O10").Value = 1500
$ws.Range("P10").Value = 1446
$ws.Range("Q10").Value = "$/envase 1 kilo"
$ws.Range("R10").Value = "Provincia de Diguillín"
$ws.Range("S10").Value = 1446
$ws.Range("T10").Value = 1
